$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.031.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.831.22'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9985'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6318'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9995'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07523'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2938'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.90'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07710'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.829.81'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.994'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6702'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.03'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009612'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.078'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.042.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '226.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9986'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.157'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9994'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.10'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1407'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.539'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.497'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.122'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.071'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05377'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.858'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7444'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.139'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.06%  '
$ws.Range('E36').Value = '  +1.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.243.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.764'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01791'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.640'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9031'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.979.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000123'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5108'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4071'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.990'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.63%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.648'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05771'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.31%  '
